$d = $word.ActiveDocument

# Locate the paragraph that begins the "所以但引壽量" question (the first
# 【問】 paragraph in the body of the document) so a brand-new Q&A paragraph
# can be inserted immediately before it.
$finder = $d.Content.Find
$finder.ClearFormatting()
$finder.Text = "所以但引壽量，不引他部者"
$finder.Execute() | Out-Null

if (-not $finder.Found) {
    throw "Could not locate anchor paragraph for insertion"
}

$anchorStart = $finder.Parent.Start
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $anchorStart -and $p.Range.End -gt $anchorStart) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index"
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a new, empty paragraph immediately before the anchor paragraph.
# Word copies the anchor paragraph's paragraph formatting (list numbering,
# style, indentation, paragraph-mark run properties) onto the new paragraph.
$anchorPara.Range.InsertParagraphBefore()

# The newly created paragraph now sits at the same index the anchor used to
# occupy (the anchor itself shifted down by one).
$newPara = $d.Paragraphs.Item($anchorIndex)
$newRange = $newPara.Range
$newRange.Collapse(1)

$run1 = "【問】『"
$run2 = "云何麤？狂華無果 ，或一華多果，或多華一果，或一華一果，或前果後華，或前華後果。初喻外道，空修梵行，無所剋獲；次喻凡夫，供養父母，報在梵天；次喻聲聞，種種苦行，止得涅槃；次喻緣覺，一遠離行，亦得涅槃；次喻須陀洹，却後修道；次喻菩薩，先籍緣修，生後真修。皆是麤華，不以為喻。"
$run3 = "』。師父，請問什麼是「一遠離行」、「却後修道」、「先籍緣修，生後真修」？【答】缘觉乐于寂静，独自修行，独自觉悟，谓远离行；须陀洹是小乘见道位（断见惑），见道之后尚须修道（断思惑），方证阿罗汉果，谓之却后修道；别教初地之前需借缘而修，初地证法身，无修而修，谓之真修。【問】须陀洹是小乘见道位（断见惑），见道之后尚须修道（断思惑），方证阿罗汉果，谓之却后修道…這是指「重慮緣真」嗎？【答】是"

$fullText = $run1 + $run2 + $run3
$newRange.InsertAfter($fullText)

$paraStart = $newPara.Range.Start

$len1 = $run1.Length
$len2 = $run2.Length
$len3 = $run3.Length

# run2 becomes bold + underlined, matching the other Q&A paragraphs' quoted
# headline text.
$run2Range = $d.Range($paraStart + $len1, $paraStart + $len1 + $len2)
$run2Range.Font.Bold = 1
$run2Range.Font.Underline = 1

Write-Host "Inserted new paragraph at index" $anchorIndex "with length" $newPara.Range.Text.Length

# The paragraph discussing "能於黑色，通達一切..." was previously split across a
# page boundary (mid-run) purely as a pagination side effect. Re-flowing the
# document after the insertion above means that stale split no longer lines
# up with anything meaningful, so collapse it back into a single run. The
# unique substring below straddles the exact run1/run2 boundary; re-typing
# just the one character after the boundary merges the two runs back
# together while preserving the surrounding run formatting/properties.
$finder2 = $d.Content.Find
$finder2.ClearFormatting()
$finder2.Text = "一切是；通达一切非非"
$finder2.Execute() | Out-Null

if ($finder2.Found) {
    $boundary = $finder2.Parent.Start + 4
    $charRange = $d.Range($boundary, $boundary + 1)
    $ch = $charRange.Text
    $charRange.Delete()
    $reinsertRange = $d.Range($boundary, $boundary)
    $reinsertRange.InsertAfter($ch)
    Write-Host "Merged split run at boundary" $boundary
} else {
    Write-Host "WARNING: merge boundary text not found"
}

